# Insert 3 new data rows just above the current row 235. This pushes the
# existing rows 235..314 down to 238..317 (dimension becomes A1:R317).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(235).EntireRow.Insert()
$ws.Rows.Item(235).EntireRow.Insert()
$ws.Rows.Item(235).EntireRow.Insert()

# --- New row 235 ---
$ws.Cells.Item(235, 1).Value = 4
$ws.Cells.Item(235, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(235, 3).Value = "Los Lagos"
$ws.Cells.Item(235, 4).Value = 44463
$ws.Cells.Item(235, 5).Value = 10
$ws.Cells.Item(235, 6).Value = 100112033
$ws.Cells.Item(235, 7).Value = "Lechuga"
$ws.Cells.Item(235, 8).Value = "Conconina(o)"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 120
$ws.Cells.Item(235, 11).Value = 9000
$ws.Cells.Item(235, 12).Value = 9000
$ws.Cells.Item(235, 13).Value = 9000
$ws.Cells.Item(235, 14).Value = "`$/caja 10 unidades"
$ws.Cells.Item(235, 15).Value = "Región Metropolitana"
$ws.Cells.Item(235, 16).Value = 900
$ws.Cells.Item(235, 17).Value = 10
$ws.Cells.Item(235, 18).Value = "Hortaliza"

# --- New row 236 ---
$ws.Cells.Item(236, 1).Value = 4
$ws.Cells.Item(236, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(236, 3).Value = "Los Lagos"
$ws.Cells.Item(236, 4).Value = 44463
$ws.Cells.Item(236, 5).Value = 10
$ws.Cells.Item(236, 6).Value = 100112033
$ws.Cells.Item(236, 7).Value = "Lechuga"
$ws.Cells.Item(236, 8).Value = "Escarola"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 450
$ws.Cells.Item(236, 11).Value = 12000
$ws.Cells.Item(236, 12).Value = 12000
$ws.Cells.Item(236, 13).Value = 12000
$ws.Cells.Item(236, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(236, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(236, 16).Value = 800
$ws.Cells.Item(236, 17).Value = 15
$ws.Cells.Item(236, 18).Value = "Hortaliza"

# --- New row 237 ---
$ws.Cells.Item(237, 1).Value = 4
$ws.Cells.Item(237, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(237, 3).Value = "Los Lagos"
$ws.Cells.Item(237, 4).Value = 44463
$ws.Cells.Item(237, 5).Value = 10
$ws.Cells.Item(237, 6).Value = 100112033
$ws.Cells.Item(237, 7).Value = "Lechuga"
$ws.Cells.Item(237, 8).Value = "Marina"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 300
$ws.Cells.Item(237, 11).Value = 9000
$ws.Cells.Item(237, 12).Value = 9000
$ws.Cells.Item(237, 13).Value = 9000
$ws.Cells.Item(237, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(237, 15).Value = "Región Metropolitana"
$ws.Cells.Item(237, 16).Value = 600
$ws.Cells.Item(237, 17).Value = 15
$ws.Cells.Item(237, 18).Value = "Hortaliza"
